$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting the existing records
# (rows 24:112) down to rows 25:113. This also naturally "duplicates"
# the former last row (old 112) into the new last row (113).
$ws.Range("A24").EntireRow.Insert()

# Populate the freshly-inserted row 24 with the new weekly record.
$ws.Range("A24").Value = 4
$ws.Range("B24").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C24").Value = "Los Lagos"
$ws.Range("D24").Value = 44487
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 100112009
$ws.Range("G24").Value = "Acelga"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 3500
$ws.Range("L24").Value = 3500
$ws.Range("M24").Value = 3500
$ws.Range("N24").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 875
$ws.Range("Q24").Value = 4
$ws.Range("R24").Value = "Hortaliza"
